# Remove client C1002 ("Tatiana Avila") — row 4 of the "Clientes" sheet.
# Deleting the entire row shifts subsequent rows (5,6) up to (4,5) and
# shrinks the used range / sheet dimension accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")
$ws.Rows(4).EntireRow.Delete()
